# Update the placeholder/reference tags in column D of the HMIS 105:2.1-2.7
# report template so that they point at the correct datasets.
#
# The "A" (Antenatal / 2.1-2.4) section tags of the form "#A...#" become
# "#A.A...#" and the "P" (Postnatal / 2.2) section tags of the form
# "#P...#" become "#P.P...#".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2.1 Antenatal / 2.4 EID section ---------------------------------
$ws.Range("D7").Value  = "#A.A1-01#"
$ws.Range("D8").Value  = "#A.A1-02#"
$ws.Range("D9").Value  = "#A.A1-03#"
$ws.Range("D10").Value = "#A.A1-04#"
$ws.Range("D11").Value = "#A.A2-01#"
$ws.Range("D12").Value = "#A.A2-02#"
$ws.Range("D13").Value = "#A.A2-03#"
$ws.Range("D14").Value = "#A.A3#"
$ws.Range("D15").Value = "#A.A4-01#"
$ws.Range("D16").Value = "#A.A4-02#"
$ws.Range("D17").Value = "#A.A4-03#"
$ws.Range("D18").Value = "#A.A5-T#"
$ws.Range("D19").Value = "#A.A5-CS#"
$ws.Range("D20").Value = "#A.A6-T#"
$ws.Range("D21").Value = "#A.A6-FSG#"
$ws.Range("D22").Value = "#A.A7-01#"
$ws.Range("D23").Value = "#A.A7-02#"
$ws.Range("D24").Value = "#A.A7-03#"
$ws.Range("D25").Value = "#A.A8-01#"
$ws.Range("D26").Value = "#A.A8-02#"
$ws.Range("D27").Value = "#A.A8-03#"
$ws.Range("D28").Value = "#A.A9#"
$ws.Range("D29").Value = "#A.A10#"
$ws.Range("D30").Value = "#A.A11#"
$ws.Range("D31").Value = "#A.A12#"
$ws.Range("D32").Value = "#A.A13-01#"
$ws.Range("D33").Value = "#A.A13-02#"
$ws.Range("D34").Value = "#A.A13-03#"
$ws.Range("D35").Value = "#A.A14-01#"
$ws.Range("D36").Value = "#A.A14-02#"
$ws.Range("D37").Value = "#A.A14-03#"
$ws.Range("D38").Value = "#A.A15-CD4#"
$ws.Range("D39").Value = "#A.A15-WHO#"
$ws.Range("D40").Value = "#A.A16#"
$ws.Range("D41").Value = "#A.A17-T#"
$ws.Range("D42").Value = "#A.A17-TRRk#"
$ws.Range("D43").Value = "#A.A18#"
$ws.Range("D44").Value = "#A.A19#"
$ws.Range("D45").Value = "#A.A20#"
$ws.Range("D46").Value = "#A.A21#"
$ws.Range("D47").Value = "#A.A22-T#"
$ws.Range("D48").Value = "#A.A22-HIVp#"

# --- 2.2 Maternity / Postnatal section --------------------------------
$ws.Range("D89").Value  = "#P.P1-A-01#"
$ws.Range("D90").Value  = "#P.P1-A-02#"
$ws.Range("D91").Value  = "#P.P1-A-03#"
$ws.Range("D92").Value  = "#P.P1-6H#"
$ws.Range("D93").Value  = "#P.P1-6D#"
$ws.Range("D94").Value  = "#P.P1-6W#"
$ws.Range("D95").Value  = "#P.P1-6M#"
$ws.Range("D96").Value  = "#P.P2-1#"
$ws.Range("D97").Value  = "#P.P2-2#"
$ws.Range("D98").Value  = "#P.P3-1#"
$ws.Range("D99").Value  = "#P.P3-2#"
$ws.Range("D100").Value = "#P.P4#"
$ws.Range("D101").Value = "#P.P5#"
$ws.Range("D102").Value = "#P.P6#"
$ws.Range("D103").Value = "#P.P7#"
$ws.Range("D104").Value = "#P.P8#"
$ws.Range("D105").Value = "#P.P9#"

# D103 previously used a stray duplicate "Calibri 12" font definition
# (missing the charset attribute of the normal body font). Re-apply the
# canonical font so the cell is re-mapped onto the regular shared style
# instead of the redundant one.
$ws.Range("D103").Font.Name = "Calibri"

# Leave the cursor where it naturally landed after typing through the
# "A" section tags (matches the saved view/selection state).
[void]$ws.Range("D49").Select()
